$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price strings that can look numeric (e.g. "529.15",
# "69.306.82", "0.0000334"). Force text format before assigning so Excel
# does not silently reinterpret them as numbers / apply float rounding,
# matching the original sheet which stores these as text.

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '69.306.82'
$ws.Range('E2').Value = '  +1.49%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.898.18'
$ws.Range('E3').Value = '  +0.03%  '
$ws.Range('E4').Value = '  -0.08%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '529.15'
$ws.Range('E5').Value = '  +8.98%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '144.45'
$ws.Range('E6').Value = '  -1.04%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.611'
$ws.Range('E7').Value = '  -1.79%  '
$ws.Range('E8').Value = '  +0.07%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.718'
$ws.Range('E9').Value = '  -3.08%  '
$ws.Range('E10').Value = '  -3.12%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0000334'
$ws.Range('E11').Value = '  -5.72%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '42.08'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '4.517.78'
$ws.Range('E13').Value = '  -0.04%  '
$ws.Range('E14').Value = '  -2.26%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '3.906.57'
$ws.Range('E15').Value = '  -0.40%  '
$ws.Range('E16').Value = '  -2.17%  '
$ws.Range('E17').Value = '  +6.58%  '
$ws.Range('E18').Value = '  -1.50%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '19.75'
$ws.Range('E19').Value = '  -1.12%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '69.282.28'
$ws.Range('E20').Value = '  +1.42%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '424.24'
$ws.Range('E21').Value = '  -1.53%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '3.38'
$ws.Range('E22').Value = '  -5.58%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '14.11'
$ws.Range('E23').Value = '  -4.37%  '
$ws.Range('E24').Value = '  -1.24%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '4.02'
$ws.Range('E25').Value = '  +8.57%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '11.39'
$ws.Range('E26').Value = '  -8.81%  '
$ws.Range('E27').Value = '  -4.33%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '36.36'
$ws.Range('E28').Value = '  -2.39%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '688.78'
$ws.Range('E29').Value = '  -4.21%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '13.17'
$ws.Range('E30').Value = '  -2.21%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.126'
$ws.Range('E31').Value = '  -2.64%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '2.84'
$ws.Range('E32').Value = '  -2.76%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '68.56'
$ws.Range('E33').Value = '  +10.99%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.0₃0864'
$ws.Range('E34').Value = '  -1.16%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.427'
$ws.Range('E35').Value = '  +7.49%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '5.90'
$ws.Range('E36').Value = '  -2.42%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '39.80'
$ws.Range('E37').Value = '  -2.56%  '
$ws.Range('E38').Value = '  +2.53%  '
$ws.Range('E39').Value = '  -0.21%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '3.27'
$ws.Range('E41').Value = '  +6.37%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.0483'
$ws.Range('E42').Value = '  -2.81%  '
$ws.Range('E43').Value = '  +7.81%  '
$ws.Range('E44').Value = '  -7.31%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '3.42'
$ws.Range('E46').Value = '  -1.42%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.000279'
$ws.Range('E47').Value = '  +13.09%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.98'
$ws.Range('E48').Value = '  +6.52%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '2.756.38'
$ws.Range('E49').Value = '  +14.63%  '
$ws.Range('E50').Value = '  -6.20%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '144.72'
$ws.Range('E51').Value = '  +0.20%  '
